# Update res_line/pl_mw data for the 380 kV (220?) case:
# rows 2-25 (A = 0..23), columns B,C,D,E,G,I,J,L get new computed values.
# Columns F,H,K,M,N,O remain unchanged (all zero in this range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","G","I","J","L")

$data = @(
    @(2, "5.035494849214842","1.12659492521243","0.346528289970081","0.05903769800093972","0.002780563059324228","7.191460915097792","0.01218124864221082","1.04832652519147"),
    @(3, "5.028737139167561","1.086597679444651","0.3468041363828789","0.05925245021986281","0.002789629326258581","6.893310474225046","0.01075984553900611","1.04052333215995"),
    @(4, "5.029255004330651","1.062833105112588","0.3471178087621354","0.05939301337577252","0.002795473260195666","6.710555348881854","0.009884759872836213","1.036297381374538"),
    @(5, "5.030634816663508","1.053346058508509","0.347281871244725","0.05945248899659838","0.002797924730585573","6.636144109231253","0.00952749221420035","1.034716767028357"),
    @(6, "5.030934387539162","1.05178258295814","0.3473113010941447","0.0594624976331779","0.002798336033484394","6.623791544624567","0.009468125965881313","1.034462834889084"),
    @(7, "5.029268886817079","1.062704364143599","0.3471198746981372","0.0593938065889148","0.002795506037518888","6.70955157779045","0.009879944419509457","1.036275492574745"),
    @(8, "5.032193159116559","1.112637613546326","0.346593421798417","0.05910994167948758","0.002783631763292038","7.088584972190347","0.01169158826046157","1.04551827217071"),
    @(9, "5.075195063219894","1.216973474420683","0.3467086424010546","0.05862207281484322","0.002762531498023111","7.835053636736717","0.01522924138902937","1.068160272567127"),
    @(10, "5.129862504543667","1.297716400405534","0.3474972546098343","0.05830519372862586","0.002748341013796141","8.386485297042725","0.01782446174347285","1.087597949046966"),
    @(11, "5.159824342801301","1.335375727169719","0.3480099319279901","0.05816998191392297","0.002742165882059419","8.63823410441637","0.01900541996751315","1.097060498781332"),
    @(12, "5.171909283178593","1.349773085534309","0.348226290381902","0.05812005990303115","0.002739867468752209","8.733712921402287","0.0194527659292234","1.10073381360732"),
    @(13, "5.169273597542656","1.346666233207372","0.3481787041581299","0.05813075467027495","0.002740360700212935","8.713143026904561","0.0193564141998479","1.099938680365028"),
    @(14, "5.160803727725693","1.336557451535839","0.3480272860774392","0.05816584918291667","0.002741975990916272","8.646086155589387","0.01904222001680722","1.097360894766013"),
    @(15, "5.15571213417121","1.330383409453646","0.3479374342161776","0.05818751208642414","0.002742970599116568","8.60503155996895","0.01884978824408279","1.095793680273289"),
    @(16, "5.128007479452435","1.295274222810463","0.3474668544790376","0.05831420952234811","0.002748750184452063","8.370052523695648","0.01774729847011969","1.086992103104677"),
    @(17, "5.112320730327042","1.273975914156324","0.3472176527727413","0.05839421960250935","0.002752367312059766","8.226143366559285","0.0170711159657202","1.081752084674946"),
    @(18, "5.103777124845124","1.261812943871917","0.3470888022274465","0.05844108090324251","0.002754474180061474","8.143454240200697","0.01668222399479902","1.07879649238987"),
    @(19, "5.100966477427619","1.257709666899643","0.3470476607801203","0.05845709205271543","0.002755192071760086","8.115470953112663","0.01655055441208475","1.077805771762598"),
    @(20, "5.113940989966579","1.27623410358018","0.347242681127625","0.05838561532865238","0.002751979533403199","8.241453961809299","0.0171430929775056","1.082303850751998"),
    @(21, "5.163271422273624","1.3395229143905","0.3480711575411561","0.05815550638496392","0.00274150045831709","8.665778252081338","0.01913450187889509","1.09811560185787"),
    @(22, "5.199822681602484","1.381683412628263","0.3487421645022835","0.05801257385386371","0.002734884636950225","8.943964034343367","0.02043687695633878","1.108974813839382"),
    @(23, "5.179917874005184","1.359107564526312","0.3483721528310753","0.05808817912579745","0.002738394423648645","8.795405945241214","0.01974166520606957","1.103130688318686"),
    @(24, "5.113206991825223","1.275212921855541","0.3472313209034041","0.05838950263428266","0.002752154762909022","8.234531899230774","0.01711055263050554","1.082054219942563"),
    @(25, "5.059536363319978","1.188043034923112","0.3465541807205454","0.05874672919719925","0.002768007831943586","7.63266480485197","0.01427328100709957","1.06154648783891")
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = [double]$entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}

Write-Host "Updated $($data.Length) rows across columns $($columns -join ', ')"
